$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 : omar ---
$ws.Range("A1").Value = 1
$ws.Range("B1").Value = "omar"
$ws.Range("C1").Value = 1
$ws.Range("D1").Value = 1

# E1 / H1 hold text that Excel would otherwise auto-convert (a date-looking
# string and a leading-zero phone number). Stage them on a scratch cell that
# is explicitly formatted as Text, then copy/paste-values into place so the
# destination keeps its original (default) style while the stored cell type
# stays text.
$ws.Range("Z1").NumberFormat = "@"
$ws.Range("Z1").Value = "1997-05-31"
$ws.Range("Z1").Copy()
$ws.Range("E1").PasteSpecial(-4163)

$ws.Range("F1").Value = "omar@gmail.com"
$ws.Range("G1").Value = "7110eda4d09e062aa5e4a390b0a572ac0d2c0220"

$ws.Range("Z1").Value = "01157979606"
$ws.Range("Z1").Copy()
$ws.Range("H1").PasteSpecial(-4163)

$ws.Range("I1").Value = 6

# --- Row 2 : Mo ---
$ws.Range("A2").Value = 2
$ws.Range("B2").Value = "Mo"
$ws.Range("C2").Value = 2
$ws.Range("D2").Value = 2

$ws.Range("Z1").Value = "2000-05-31"
$ws.Range("Z1").Copy()
$ws.Range("E2").PasteSpecial(-4163)

$ws.Range("F2").Value = "mo@gmail.com"
$ws.Range("G2").Value = "ac1ab23d6288711be64a25bf13432baf1e60b2bd"

$ws.Range("Z1").Value = "012825347698"
$ws.Range("Z1").Copy()
$ws.Range("H2").PasteSpecial(-4163)

$ws.Range("I2").Value = 5

# clean up the scratch cell so it doesn't leave stray data on the sheet
$ws.Range("Z1").Clear()
